$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.973.91"
$ws.Range("E2").Value = "  +4.96%  "
$ws.Range("D3").Value = "3.514.94"
$ws.Range("E3").Value = "  +2.90%  "
$ws.Range("E4").Value = "  +0.05%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "593.23"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +4.04%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "169.00"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +6.93%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "3.517.14"
$ws.Range("E8").Value = "  +2.92%  "
$ws.Range("E9").Value = "  +1.27%  "
$ws.Range("E10").Value = "  +0.35%  "
$ws.Range("E11").Value = "  +5.61%  "
$ws.Range("E12").Value = "  +4.30%  "
$ws.Range("D13").Value = "4.123.59"
$ws.Range("E13").Value = "  +2.94%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.134"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.05%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "28.18"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +4.07%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.0000179"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +4.15%  "
$ws.Range("D17").Value = "66.935.79"
$ws.Range("E17").Value = "  +4.81%  "
$ws.Range("D18").Value = "3.525.01"
$ws.Range("E18").Value = "  +2.25%  "
$ws.Range("E19").Value = "  +4.03%  "
$ws.Range("E20").Value = "  +3.12%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "395.06"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +3.46%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "7.97"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +2.08%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "73.59"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +3.29%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "0.0000127"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +10.22%  "
$ws.Range("E25").Value = "  -0.34%  "
$ws.Range("E26").Value = "  +3.24%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "10.18"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +5.10%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.182"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +2.18%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "6.40"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +5.34%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.47"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +5.97%  "
$ws.Range("E32").Value = "  +4.24%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "23.60"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +3.17%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "7.46"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +7.34%  "
$ws.Range("E36").Value = "  +5.67%  "
$ws.Range("E37").Value = "  +0.47%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.900"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +6.62%  "
$ws.Range("E39").Value = "  +5.76%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.0752"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +4.06%  "
$ws.Range("E41").Value = "  +7.07%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "26.58"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +1.82%  "
$ws.Range("E43").Value = "  +4.96%  "
$ws.Range("D44").Value = "2.837.55"
$ws.Range("E44").Value = "  +0.93%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "43.52"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +1.05%  "
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "2.57"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +7.19%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "26.36"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.03%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.0315"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +3.56%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "352.19"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +5.76%  "
$ws.Range("E50").Value = "  +4.66%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "33.50"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +11.42%  "
